# Update the "Förändrad" (changed) date column (C) for rows 2 through 16
# from 2023-09-05 (45174) to 2023-09-06 (45175), keeping the existing
# date formatting/style intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 16; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45174) {
        $cell.Value = 45175
    }
}
